$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 36; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # column C ("Förändrad")
    if ($cell.Value2 -eq 46060) {
        $cell.Value = 46061
    }
}
